$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 3140
$ws.Range("L3").Value = 3216
$ws.Range("B4").Value = 1712
$ws.Range("J4").Value = 1868
$ws.Range("L4").Value = 818
$ws.Range("L6").Value = 2853
$ws.Range("B7").Value = 23344
$ws.Range("J7").Value = 29343
$ws.Range("L7").Value = 10206

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 87
$ws.Range("L4").Value = 38
$ws.Range("L5").Value = 41
$ws.Range("L7").Value = 340
$ws.Range("L8").Value = 649
$ws.Range("L14").Value = 49
$ws.Range("L19").Value = 289
$ws.Range("L20").Value = 258
$ws.Range("L23").Value = 108
$ws.Range("L29").Value = 555
$ws.Range("J33").Value = 1313
$ws.Range("L33").Value = 470
$ws.Range("L34").Value = 66
$ws.Range("L36").Value = 141
$ws.Range("L41").Value = 46
$ws.Range("L42").Value = 328
$ws.Range("L43").Value = 80
$ws.Range("L48").Value = 137
$ws.Range("L49").Value = 56
$ws.Range("L52").Value = 201
$ws.Range("L54").Value = 214
$ws.Range("B63").Value = 416
$ws.Range("L63").Value = 34
$ws.Range("L65").Value = 194
$ws.Range("L67").Value = 374
$ws.Range("L76").Value = 137
$ws.Range("L78").Value = 124
$ws.Range("L79").Value = 263
$ws.Range("L84").Value = 103
$ws.Range("L85").Value = 520
$ws.Range("L86").Value = 76
$ws.Range("L89").Value = 137
$ws.Range("L91").Value = 149
$ws.Range("L92").Value = 29
$ws.Range("L93").Value = 53
$ws.Range("L94").Value = 120
$ws.Range("L97").Value = 92
$ws.Range("L98").Value = 63
$ws.Range("L99").Value = 172
$ws.Range("B101").Value = 23344
$ws.Range("J101").Value = 29343
$ws.Range("L101").Value = 10206

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 104
$ws.Range("L7").Value = 340

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 42
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 154
$ws.Range("L3").Value = 214
$ws.Range("L7").Value = 520

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 59
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 191
$ws.Range("L3").Value = 213
$ws.Range("L6").Value = 179
$ws.Range("L7").Value = 649

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 133
$ws.Range("J4").Value = 57
$ws.Range("L6").Value = 162
$ws.Range("J7").Value = 1313
$ws.Range("L7").Value = 470

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 112
$ws.Range("L6").Value = 113

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 71
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 194

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 43
$ws.Range("L3").Value = 70
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 172

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 111
$ws.Range("L3").Value = 141
$ws.Range("L6").Value = 84
$ws.Range("L7").Value = 374

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 47
$ws.Range("L6").Value = 106
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 205
$ws.Range("L6").Value = 148
$ws.Range("L7").Value = 555

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 31
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 101
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 289

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 98
$ws.Range("L3").Value = 102
$ws.Range("L6").Value = 91
$ws.Range("L7").Value = 328

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 149

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 88
$ws.Range("L3").Value = 94
$ws.Range("L7").Value = 263

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 84
$ws.Range("L3").Value = 78
$ws.Range("L7").Value = 258

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 53
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 141

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 120

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L4").Value = 41
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 38
